$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 995.0769
$ws.Range("I2").Value = 93.7
$ws.Range("K2").Value = 93.7
$ws.Range("M2").Value = 19.3
$ws.Range("H19").Value = 1536.9546
$ws.Range("I19").Value = 1966.6364
$ws.Range("K19").Value = 1966.6364
$ws.Range("M19").Value = -1791.6364
$ws.Range("H113").Value = 35724148
$ws.Range("I113").Value = 76927090
$ws.Range("J113").Value = 14929.866
$ws.Range("K113").Value = 76927090
$ws.Range("L113").Value = 14929.866
$ws.Range("M113").Value = -76923836
$ws.Range("N113").Value = -21437.866
$ws.Range("H116").Value = 3499.3333
$ws.Range("I116").Value = 3200
$ws.Range("J116").Value = 3649
$ws.Range("K116").Value = 3200
$ws.Range("L116").Value = 3649
$ws.Range("M116").Value = 242
$ws.Range("N116").Value = -10533
$ws.Range("H132").Value = 1985.5555
$ws.Range("I132").Value = 1501.3226
$ws.Range("K132").Value = 4503.9678
$ws.Range("M132").Value = -1973.9678
$ws.Range("H137").Value = 2176.6667
$ws.Range("I137").Value = 2266.4
$ws.Range("J137").Value = 1920.2858
$ws.Range("K137").Value = 6799.200000000001
$ws.Range("L137").Value = 5760.857400000001
$ws.Range("M137").Value = -4249.200000000001
$ws.Range("N137").Value = -10860.8574
$ws.Range("H138").Value = 2983.5
$ws.Range("I138").Value = 1691.9333
$ws.Range("J138").Value = 4003.158
$ws.Range("K138").Value = 5075.7999
$ws.Range("L138").Value = 12009.474
$ws.Range("M138").Value = 64.20010000000002
$ws.Range("N138").Value = -22289.474
$ws.Range("H141").Value = 5616.615
$ws.Range("I141").Value = 5334.6665
$ws.Range("K141").Value = 16003.9995
$ws.Range("M141").Value = -10823.9995

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3153.3
$ws.Range("I2").Value = 2531.4443
$ws.Range("K2").Value = 2531.4443
$ws.Range("M2").Value = -2418.4443
$ws.Range("H30").Value = 3018.5
$ws.Range("J30").Value = 4999.6665
$ws.Range("L30").Value = 4999.6665
$ws.Range("N30").Value = -5299.6665
$ws.Range("H45").Value = 4935
$ws.Range("I45").Value = 4150.143
$ws.Range("J45").Value = 6155.8887
$ws.Range("K45").Value = 4150.143
$ws.Range("L45").Value = 6155.8887
$ws.Range("M45").Value = -3773.143
$ws.Range("N45").Value = -6909.8887
$ws.Range("H61").Value = 8162.619
$ws.Range("I61").Value = 8370.75
$ws.Range("K61").Value = 8370.75
$ws.Range("M61").Value = -8158.75
$ws.Range("H116").Value = 3153.3
$ws.Range("I116").Value = 2531.4443
$ws.Range("K116").Value = 2531.4443
$ws.Range("M116").Value = -237.4443000000001
$ws.Range("H136").Value = 8162.619
$ws.Range("I136").Value = 8370.75
$ws.Range("K136").Value = 25112.25
$ws.Range("M136").Value = -22562.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3153.3
$ws.Range("I3").Value = 2531.4443
$ws.Range("K3").Value = 2531.4443
$ws.Range("M3").Value = -2417.4443
$ws.Range("H22").Value = 3937.4
$ws.Range("I22").Value = 4609.25
$ws.Range("K22").Value = 4609.25
$ws.Range("M22").Value = -4436.25
$ws.Range("H43").Value = 262450
$ws.Range("J43").Value = 262450
$ws.Range("L43").Value = 262450
$ws.Range("N43").Value = -262812
$ws.Range("H105").Value = 1739.9333
$ws.Range("I105").Value = 2229.75
$ws.Range("K105").Value = 2229.75
$ws.Range("M105").Value = -482.75
$ws.Range("H107").Value = 2548.5625
$ws.Range("I107").Value = 2417.5833
$ws.Range("K107").Value = 2417.5833
$ws.Range("M107").Value = -497.5832999999998
$ws.Range("H134").Value = 9314.526
$ws.Range("I134").Value = 9233.883
$ws.Range("K134").Value = 27701.649
$ws.Range("M134").Value = -25166.649

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3967.05
$ws.Range("I16").Value = 2614.923
$ws.Range("J16").Value = 6478.143
$ws.Range("K16").Value = 2614.923
$ws.Range("L16").Value = 6478.143
$ws.Range("M16").Value = -2327.923
$ws.Range("N16").Value = -7052.143
$ws.Range("H68").Value = 41400
$ws.Range("I68").Value = 42800
$ws.Range("K68").Value = 42800
$ws.Range("M68").Value = -42051
$ws.Range("H71").Value = 41400
$ws.Range("I71").Value = 42800
$ws.Range("K71").Value = 128400
$ws.Range("M71").Value = -124656
$ws.Range("H74").Value = 40716.332
$ws.Range("J74").Value = 40716.332
$ws.Range("L74").Value = 40716.332
$ws.Range("N74").Value = -42464.332
$ws.Range("H77").Value = 40716.332
$ws.Range("J77").Value = 40716.332
$ws.Range("L77").Value = 122148.996
$ws.Range("N77").Value = -130884.996
$ws.Range("H97").Value = 26016.857
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H107").Value = 2167.7144
$ws.Range("I107").Value = 637
$ws.Range("K107").Value = 637
$ws.Range("M107").Value = 1283
$ws.Range("H113").Value = 3967.05
$ws.Range("I113").Value = 2614.923
$ws.Range("J113").Value = 6478.143
$ws.Range("K113").Value = 2614.923
$ws.Range("L113").Value = 6478.143
$ws.Range("M113").Value = -444.9229999999998
$ws.Range("N113").Value = -10818.143
$ws.Range("H132").Value = 4339.75
$ws.Range("I132").Value = 2442.4546
$ws.Range("K132").Value = 7327.3638
$ws.Range("M132").Value = -4797.3638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 389.7143
$ws.Range("I60").Value = 479.6
$ws.Range("K60").Value = 1438.8
$ws.Range("M60").Value = -1187.8
$ws.Range("H107").Value = 323.18182
$ws.Range("I107").Value = 323.18182
$ws.Range("K107").Value = 969.54546
$ws.Range("M107").Value = 950.45454

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 9091.583000000001
$ws.Range("I70").Value = 6871.353
$ws.Range("K70").Value = 6871.353
$ws.Range("M70").Value = -6601.353
$ws.Range("H73").Value = 9091.583000000001
$ws.Range("I73").Value = 6871.353
$ws.Range("K73").Value = 6871.353
$ws.Range("M73").Value = -5935.353
$ws.Range("H102").Value = 6077.25
$ws.Range("I102").Value = 4116
$ws.Range("K102").Value = 4116
$ws.Range("M102").Value = -2494
$ws.Range("H107").Value = 794227.2
$ws.Range("I107").Value = 1304130.6
$ws.Range("J107").Value = 1044.1111
$ws.Range("K107").Value = 1304130.6
$ws.Range("L107").Value = 1044.1111
$ws.Range("M107").Value = -1302210.6
$ws.Range("N107").Value = -4884.1111
$ws.Range("H113").Value = 4292.4287
$ws.Range("J113").Value = 4599.6665
$ws.Range("L113").Value = 4599.6665
$ws.Range("N113").Value = -8939.666499999999
$ws.Range("H122").Value = 2021.6428
$ws.Range("I122").Value = 1835.3
$ws.Range("K122").Value = 5505.9
$ws.Range("M122").Value = -3055.9
$ws.Range("H132").Value = 7618.6
$ws.Range("I132").Value = 7020.125
$ws.Range("K132").Value = 21060.375
$ws.Range("M132").Value = -18530.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 109999.164
$ws.Range("J6").Value = 109999.164
$ws.Range("L6").Value = 109999.164
$ws.Range("N6").Value = -110223.164
$ws.Range("H22").Value = 5322
$ws.Range("I22").Value = 5899
$ws.Range("J22").Value = 5249.875
$ws.Range("K22").Value = 5899
$ws.Range("L22").Value = 5249.875
$ws.Range("M22").Value = -5604
$ws.Range("N22").Value = -5839.875
$ws.Range("H27").Value = 5322
$ws.Range("I27").Value = 5899
$ws.Range("J27").Value = 5249.875
$ws.Range("K27").Value = 5899
$ws.Range("L27").Value = 5249.875
$ws.Range("M27").Value = -5792
$ws.Range("N27").Value = -5463.875
$ws.Range("I96").Value = 40000
$ws.Range("K96").Value = 40000
$ws.Range("M96").Value = -37254
$ws.Range("H122").Value = 2988
$ws.Range("I122").Value = 2988
$ws.Range("K122").Value = 8964
$ws.Range("M122").Value = -6514
$ws.Range("H136").Value = 6297.5093
$ws.Range("I136").Value = 6276.7114
$ws.Range("K136").Value = 18830.1342
$ws.Range("M136").Value = -16280.1342

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H49").Value = 9933.333000000001
$ws.Range("I49").Value = 9900
$ws.Range("J49").Value = 10000
$ws.Range("K49").Value = 9900
$ws.Range("L49").Value = 10000
$ws.Range("M49").Value = -9670
$ws.Range("N49").Value = -10460
$ws.Range("H100").Value = 632.9524
$ws.Range("I100").Value = 535.5625
$ws.Range("K100").Value = 1071.125
$ws.Range("M100").Value = -530.125
$ws.Range("H107").Value = 823.5769
$ws.Range("I107").Value = 828.0625
$ws.Range("J107").Value = 816.4
$ws.Range("K107").Value = 2484.1875
$ws.Range("L107").Value = 2449.2
$ws.Range("M107").Value = -564.1875
$ws.Range("N107").Value = -6289.2
$ws.Range("H113").Value = 1133.05
$ws.Range("J113").Value = 2799.6
$ws.Range("L113").Value = 8398.799999999999
$ws.Range("N113").Value = -12738.8
$ws.Range("H136").Value = 2670.9
$ws.Range("I136").Value = 2022
$ws.Range("K136").Value = 6066
$ws.Range("M136").Value = -3516
